# "Added support for maces" - adds a new weapon-class (game_class_id = 13)
# block of Class Specials rows, tweaks the "grow by" wording/values on the
# existing gun specials (rows 112-117), and widens column C to fit the new,
# longer specialty names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class Specials")

# ---------------------------------------------------------------------------
# 1. Fix the "grow by" amounts for the existing gun specials (rows 112-117):
#    the per-level growth amount was 10x too large, and the description text
#    needs the same number corrected.
# ---------------------------------------------------------------------------

$ws.Range("G112").Value = 50
$ws.Range("D112").Value = "Deal 5,000 Damage using 20% of your damage stat as bonus - which will then grow by 50 for an additional 5,000 damage - only while using ATTACK. You will also grow your damage by +150%"

$ws.Range("G113").Value = 100
$ws.Range("D113").Value = "Deal 10,000 Damage growing by 100 for an additional 10,000 damage while applying 25% of your damage stat to the damage - only while using ATTACK. Reduce all aspects of the enemy by 100% over time."

$ws.Range("G114").Value = 200
$ws.Range("D114").Value = "Deal 20,000 Damage while growing it by 200 for an additional 20,000 and applying 27% of your damage stat as bonus damage. Damage is only dealt while using ATTACK AND CAST. You will also grow your spell damage by 125% over time and reduce the enemies spell evasion and affix damage by 100% over time."

$ws.Range("G115").Value = 250
$ws.Range("D115").Value = "Deal 25,000 Damage growing by 250 damage for an additional 25,000 damage with 30% of your damage stat as bonus. This only procs while using ATTACK. You will also grow your own damage by +250% over time."

$ws.Range("G116").Value = 500
$ws.Range("D116").Value = "Deal 50,000 in damage growing by 500 for an additional 50,000 damage while applying 12% of your damage stat as bonus damage. Damage only procs during Cast and Attack. Reduce ass aspects on an enemy by 100% over time and grow your own spell damage by +200%"

$ws.Range("G117").Value = 750
$ws.Range("D117").Value = "Deal 75,000 Damage growing by 750 damage while applying 40% of your damage stat to the over all damage as a bonus. Only procs when using DEFEND. Will also grow your base damage stat by +300% overtime."

# ---------------------------------------------------------------------------
# 2. Append the new mace (game_class_id = 13) specials as rows 118-126.
# ---------------------------------------------------------------------------

# Row 118 - Clerical Prayer
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = 13
$ws.Range("C118").Value = "Clerical Prayer"
$ws.Range("D118").Value = "Over time grow your health by +50% and boost your healing - done through healing spells - by +150% over time."
$ws.Range("E118").Value = 1
$ws.Range("K118").Value = 0.015
$ws.Range("M118").Value = 0.005

# Row 119 - Churches Grace
$ws.Range("A119").Value = 118
$ws.Range("B119").Value = 13
$ws.Range("C119").Value = "Churches Grace"
$ws.Range("D119").Value = "Increase Damage, Healing and Armour class by 150% and 170% (For AC overtime. Increase your health and damage stat by 75% over time,"
$ws.Range("E119").Value = 12
$ws.Range("I119").Value = 0.015
$ws.Range("J119").Value = 0.017
$ws.Range("K119").Value = 0.015
$ws.Range("M119").Value = 0.0075
$ws.Range("N119").Value = 0.0075

# Row 120 - Blessed Rage
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = 13
$ws.Range("C120").Value = "Blessed Rage"
$ws.Range("D120").Value = "Reduce all aspects of the enemy (see reductions) by 100% over time. Increase your Armour class by 100%, Base Healing (through healing spells) by 170%, Spell Damage by 50%, health and damage stat by 100% - all overtime."
$ws.Range("E120").Value = 24
$ws.Range("J120").Value = 0.01
$ws.Range("K120").Value = 0.017
$ws.Range("L120").Value = 0.005
$ws.Range("M120").Value = 0.01
$ws.Range("N120").Value = 0.01
$ws.Range("P120").Value = 0.01
$ws.Range("Q120").Value = 0.01
$ws.Range("R120").Value = 0.01
$ws.Range("S120").Value = 0.01
$ws.Range("T120").Value = 0.01

# Row 121 - Faithless War Cry
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = 13
$ws.Range("C121").Value = "Faithless War Cry"
$ws.Range("D121").Value = "Deal 5,000 damage growing by 50 for an additional 5,000 damage while applying 5% of your damage stat as bonus damage. must use ATTACK for this to proc. You will grow your damage stat by 150% over time."
$ws.Range("E121").Value = 36
$ws.Range("F121").Value = 5000
$ws.Range("G121").Value = 50
$ws.Range("H121").Value = 0.05
$ws.Range("N121").Value = 0.015
$ws.Range("O121").Value = "attack"

# Row 122 - Malicious Prayer
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = 13
$ws.Range("C122").Value = "Malicious Prayer"
$ws.Range("D122").Value = "Deal 10,000 damage growing by 100 over time dealing an additional 10,000 damage, while applying 10% of your damage stat towards the damage. This will only proc while using ATTACK AND CAST. Grow your healing and damage stat by 200% over time and grow your health by 100%, overtime."
$ws.Range("E122").Value = 48
$ws.Range("F122").Value = 10000
$ws.Range("G122").Value = 100
$ws.Range("H122").Value = 0.1
$ws.Range("K122").Value = 0.02
$ws.Range("M122").Value = 0.01
$ws.Range("N122").Value = 0.02
$ws.Range("O122").Value = "attack_and_cast"

# Row 123 - The Churches Holy Magic
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = 13
$ws.Range("C123").Value = "The Churches Holy Magic"
$ws.Range("D123").Value = "Grow healing (Healing Spells) by 250%, Spell damage by 150% and your damage stat by 200% over time. Deal 25,000 damage, growing by 250 damage over time for an additional 25,000 damage., which will apply 20% of your damage stat as bonus damage - however this only procs if you use: CAST"
$ws.Range("E123").Value = 60
$ws.Range("F123").Value = 25000
$ws.Range("G123").Value = 250
$ws.Range("H123").Value = 0.2
$ws.Range("K123").Value = 0.025
$ws.Range("L123").Value = 0.015
$ws.Range("N123").Value = 0.02
$ws.Range("O123").Value = "cast"

# Row 124 - Churches Defence
$ws.Range("A124").Value = 123
$ws.Range("B124").Value = 13
$ws.Range("C124").Value = "Churches Defence"
$ws.Range("D124").Value = "Deal 40,000 damage, growing by 400 for an additional, 40,000 damage while applying 18% of your damage stat as bonus damage. Only procs if you use DEFEND. Grow your own AC by 200% over time."
$ws.Range("E124").Value = 70
$ws.Range("F124").Value = 40000
$ws.Range("G124").Value = 400
$ws.Range("H124").Value = 0.18
$ws.Range("J124").Value = 0.02
$ws.Range("O124").Value = "defend"

# Row 125 - Wrath of the true God
$ws.Range("A125").Value = 124
$ws.Range("B125").Value = 13
$ws.Range("C125").Value = "Wrath of the true God"
$ws.Range("D125").Value = "Deal 60,000 Damage growing by 600 for an additional 60,000 damage while applying 15% of your damage stat as bonus damage. Only procs during ATTACK AND CAST. Grow your damage, ac and healing by 200% over time. Grow your damage stat by 300% over time."
$ws.Range("E125").Value = 80
$ws.Range("F125").Value = 60000
$ws.Range("G125").Value = 600
$ws.Range("H125").Value = 0.15
$ws.Range("I125").Value = 0.02
$ws.Range("J125").Value = 0.02
$ws.Range("K125").Value = 0.02
$ws.Range("N125").Value = 0.03
$ws.Range("O125").Value = "attack_and_cast"
$ws.Range("P125").Value = 0.01
$ws.Range("Q125").Value = 0.01
$ws.Range("R125").Value = 0.01
$ws.Range("S125").Value = 0.01
$ws.Range("T125").Value = 0.01

# Row 126 - Churches Blessing on the Faithless
$ws.Range("A126").Value = 125
$ws.Range("B126").Value = 13
$ws.Range("C126").Value = "Churches Blessing on the Faithless"
$ws.Range("D126").Value = "Deal 75,000 damage growing by 750 over time while adding 10% of your damage stat as bonus damage. Will proc during ANY attack type. Will grow all your modifiers by 200% overtime - such as Damage, Healing, Spell Damage, Health and Damage Stat."
$ws.Range("E126").Value = 90
$ws.Range("F126").Value = 75000
$ws.Range("G126").Value = 750
$ws.Range("H126").Value = 0.1
$ws.Range("I126").Value = 0.02
$ws.Range("J126").Value = 0.02
$ws.Range("K126").Value = 0.02
$ws.Range("L126").Value = 0.02
$ws.Range("M126").Value = 0.02
$ws.Range("N126").Value = 0.02
$ws.Range("O126").Value = "any"

# ---------------------------------------------------------------------------
# 3. Column C ("name") now holds longer specialty names (e.g. "Churches
#    Blessing on the Faithless"), so its best-fit width grows too.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 40.3

Write-Output "Added mace (game_class_id 13) specials rows 118-126 and updated gun special growth values."
